# Automatic update of files.
#
# - Column C ("Förändrad") bumped from 45183 to 45184 for every data row (2..10).
# - For rows 2-4 (the only rows that carry link formulas), the HYPERLINK()
#   formulas in columns S, T, V, W, X gain a second ("friendly text") argument,
#   and column Y (previously a literal/inline-string, not a real formula) is
#   rewritten as an actual HYPERLINK formula with a comma-separated second
#   argument as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C: "Förändrad" date bump (45183 -> 45184) for rows 2 through 10 ---
for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 3).Value = 45184
}

# --- Rows 2-4: update link formulas with the "ID" friendly-name argument ---
$baseUrl = "https://klasma.github.io/LoggingDetectiveFiles/Logging_OSTERSUND"

$ids = @{
    2 = "A 30683-2023"
    3 = "A 32699-2023"
    4 = "A 29992-2023"
}

foreach ($row in 2..4) {
    $id = $ids[$row]

    # S: artfynd link - the replacement text lands *inside* the URL string
    # literal (matches the source edit exactly, including the malformed quoting).
    $ws.Range("S$row").Formula = '=HYPERLINK("' + $baseUrl + '/artfynd/' + $id + '.xlsx, "' + $id + '"")'

    # T: kartor link
    $ws.Range("T$row").Formula = '=HYPERLINK("' + $baseUrl + '/kartor/' + $id + '.png", "' + $id + '")'

    # V: klagomål link
    $ws.Range("V$row").Formula = '=HYPERLINK("' + $baseUrl + '/klagomål/' + $id + '.docx", "' + $id + '")'

    # W: klagomålsmail link
    $ws.Range("W$row").Formula = '=HYPERLINK("' + $baseUrl + '/klagomålsmail/' + $id + '.docx", "' + $id + '")'

    # X: tillsyn link
    $ws.Range("X$row").Formula = '=HYPERLINK("' + $baseUrl + '/tillsyn/' + $id + '.docx", "' + $id + '")'

    # Y: tillsynsmail link - was stored as inline text, now becomes a real formula
    $ws.Range("Y$row").Formula = '=HYPERLINK("' + $baseUrl + '/tillsynsmail/' + $id + '.docx", "' + $id + '")'
}
